$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.561.70'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '3.149.11'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '563.80'
$ws.Range("E5").Value = '  +2.61%  '
$ws.Range("D6").Value = '142.90'
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.138.87'
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("D10").Value = '6.80'
$ws.Range("E10").Value = '  +5.54%  '
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").Value = '0.466'
$ws.Range("E12").Value = '  +2.46%  '
$ws.Range("D13").Value = '36.82'
$ws.Range("E13").Value = '  +3.52%  '
$ws.Range("D14").Value = '0.0000222'
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").Value = '3.651.40'
$ws.Range("E15").Value = '  +3.00%  '
$ws.Range("D16").Value = '64.628.72'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").Value = '3.147.70'
$ws.Range("E18").Value = '  +2.96%  '
$ws.Range("D19").Value = '516.41'
$ws.Range("E19").Value = '  +6.67%  '
$ws.Range("D20").Value = '6.84'
$ws.Range("E20").Value = '  +4.26%  '
$ws.Range("D21").Value = '14.02'
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").Value = '0.716'
$ws.Range("E22").Value = '  +4.95%  '
$ws.Range("E23").Value = '  +4.51%  '
$ws.Range("D24").Value = '12.77'
$ws.Range("E24").Value = '  +3.87%  '
$ws.Range("D25").Value = '79.06'
$ws.Range("E25").Value = '  +1.56%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  +16.06%  '
$ws.Range("D28").Value = '2.83'
$ws.Range("E28").Value = '  +4.92%  '
$ws.Range("D29").Value = '2.15'
$ws.Range("E29").Value = '  +4.31%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '26.63'
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("D32").Value = '2.60'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D34").Value = '550.84'
$ws.Range("E34").Value = '  -5.65%  '
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").Value = '6.09'
$ws.Range("E36").Value = '  +3.60%  '
$ws.Range("D37").Value = '53.91'
$ws.Range("E37").Value = '  +4.16%  '
$ws.Range("D38").Value = '0.0434'
$ws.Range("E38").Value = '  +6.68%  '
$ws.Range("D39").Value = '0.0824'
$ws.Range("E39").Value = '  +4.59%  '
$ws.Range("D40").Value = '3.160.17'
$ws.Range("E40").Value = '  +8.24%  '
$ws.Range("E41").Value = '  +3.80%  '
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("D43").Value = '8.29'
$ws.Range("E43").Value = '  +1.51%  '
$ws.Range("D44").Value = '0.265'
$ws.Range("E44").Value = '  +10.21%  '
$ws.Range("E45").Value = '  +8.14%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '25.35'
$ws.Range("E47").Value = '  +3.26%  '
$ws.Range("D48").Value = '120.61'
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("D49").Value = '0.108'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").Value = '0.0₃0518'
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = '2.10'
$ws.Range("E51").Value = '  +2.87%  '
